$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Ggupta"
$ws.Range("B3").Value = "Gopesh Gupta"

$ws.Range("F9").Select()
